# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.63 = 18316.85 pesos`n✅ 18316.85 pesos = 4.59 = 940.82 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update N10, O10, N12, O12 values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 216
$ws2.Range("O10").Value = 3956.44
$ws2.Range("N12").Value = 3991.15
$ws2.Range("O12").Value = 205
